$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the header row (row 1) entirely, shifting all data rows up by one.
$ws.Rows.Item(1).Delete()

# Leave the selection on the (now first) row, matching the selection Excel
# leaves behind after an entire-row delete.
$ws.Rows.Item(1).Select() | Out-Null
